$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "19.934.10"
$ws.Range("E2").Value = "  -7.48%  "

$ws.Range("D3").Value = "1.407.56"
$ws.Range("E3").Value = "  -8.07%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").Value = "'275.19"
$ws.Range("E6").Value = "  -4.70%  "

$ws.Range("D7").Value = "'0.3668"
$ws.Range("E7").Value = "  -5.23%  "

$ws.Range("D8").Value = "'0.3109"
$ws.Range("E8").Value = "  -1.93%  "

$ws.Range("D9").Value = "'39.75"
$ws.Range("E9").Value = "  -6.91%  "

$ws.Range("D10").Value = "'1.027"
$ws.Range("E10").Value = "  -3.72%  "

$ws.Range("D11").Value = "'0.06482"
$ws.Range("E11").Value = "  -9.42%  "

$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").Value = "'5.468"
$ws.Range("E13").Value = "  -4.57%  "

$ws.Range("D14").Value = "'17.61"
$ws.Range("E14").Value = "  -3.03%  "

$ws.Range("D15").Value = "'6.168"
$ws.Range("E15").Value = "  -5.86%  "

$ws.Range("D16").Value = "1.412.06"
$ws.Range("E16").Value = "  -8.34%  "

$ws.Range("D17").Value = "'0.00001017"
$ws.Range("E17").Value = "  -6.43%  "

$ws.Range("D18").Value = "'0.05695"
$ws.Range("E18").Value = "  -13.95%  "

$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").Value = "'70.82"
$ws.Range("E20").Value = "  -15.27%  "

$ws.Range("D21").Value = "'5.614"
$ws.Range("E21").Value = "  -7.95%  "

$ws.Range("E22").Value = "  -4.66%  "

$ws.Range("D23").Value = "'11.01"
$ws.Range("E23").Value = "  +1.95%  "

$ws.Range("D24").Value = "'2.257"
$ws.Range("E24").Value = "  -4.80%  "

$ws.Range("D25").Value = "19.949.41"
$ws.Range("E25").Value = "  -7.38%  "

$ws.Range("E26").Value = "  -5.37%  "

$ws.Range("D27").Value = "'133.51"
$ws.Range("E27").Value = "  -10.86%  "

$ws.Range("D28").Value = "'17.05"
$ws.Range("E28").Value = "  -7.04%  "

$ws.Range("D29").Value = "1.570.08"
$ws.Range("E29").Value = "  -8.19%  "

$ws.Range("D30").Value = "'109.36"
$ws.Range("E30").Value = "  -6.21%  "

$ws.Range("D31").Value = "'3.982"
$ws.Range("E31").Value = "  -17.70%  "

$ws.Range("D32").Value = "'5.307"
$ws.Range("E32").Value = "  -12.18%  "

$ws.Range("D33").Value = "'0.8205"
$ws.Range("E33").Value = "  -13.60%  "

$ws.Range("E34").Value = "  -3.85%  "

$ws.Range("D35").Value = "'8.441"
$ws.Range("E35").Value = "  -0.96%  "

$ws.Range("D36").Value = "'1.484"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").Value = "'0.05877"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").Value = "'4.872"
$ws.Range("E38").Value = "  -5.60%  "

$ws.Range("D39").Value = "'1.003"
$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").Value = "'0.02066"
$ws.Range("E40").Value = "  -6.30%  "

$ws.Range("D41").Value = "'10.47"
$ws.Range("E41").Value = "  -7.22%  "

$ws.Range("D42").Value = "'0.1900"
$ws.Range("E42").Value = "  -6.01%  "

$ws.Range("D43").Value = "'1.088"
$ws.Range("E43").Value = "  -7.48%  "

$ws.Range("D44").Value = "'12.37"
$ws.Range("E44").Value = "  -5.61%  "

$ws.Range("D45").Value = "'0.5291"
$ws.Range("E45").Value = "  -7.94%  "

$ws.Range("D46").Value = "'3.531"
$ws.Range("E46").Value = "  -5.05%  "

$ws.Range("D47").Value = "'0.5164"
$ws.Range("E47").Value = "  -6.91%  "

$ws.Range("D48").Value = "'114.83"
$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("D49").Value = "'1.765"
$ws.Range("E49").Value = "  -6.46%  "

$ws.Range("D50").Value = "'1.038"
$ws.Range("E50").Value = "  -10.32%  "

$ws.Range("E51").Value = "  +0.15%  "
